# Refresh the cryptos list (price + 1h volume change) as of the
# "Mon Aug 12 03:17:32 UTC 2024" GitHub Actions scrape.
# D-column values that would otherwise be auto-parsed as numbers by Excel
# (e.g. "507.10", "1.00") are forced to Text format first so the literal
# string (with trailing zeros / decimal formatting) round-trips exactly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.505.34"
$ws.Range("E2").Value = "  -4.20%  "
$ws.Range("D3").Value = "2.532.16"
$ws.Range("E3").Value = "  -3.65%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "507.10"
$ws.Range("E5").Value = "  -4.97%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.77"
$ws.Range("E6").Value = "  -7.57%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.18%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.561"
$ws.Range("E8").Value = "  -5.22%  "
$ws.Range("D9").Value = "2.534.92"
$ws.Range("E9").Value = "  -3.80%  "
$ws.Range("E10").Value = "  -7.62%  "
$ws.Range("E11").Value = "  -7.20%  "
$ws.Range("E12").Value = "  -5.12%  "
$ws.Range("E13").Value = "  -0.48%  "
$ws.Range("D14").Value = "2.977.76"
$ws.Range("E14").Value = "  -3.49%  "
$ws.Range("D15").Value = "58.484.12"
$ws.Range("E15").Value = "  -4.18%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.62"
$ws.Range("E16").Value = "  -5.30%  "
$ws.Range("E17").Value = "  -6.74%  "
$ws.Range("D18").Value = "2.539.55"
$ws.Range("E18").Value = "  -3.31%  "
$ws.Range("E19").Value = "  -5.50%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "334.70"
$ws.Range("E20").Value = "  -6.03%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.05"
$ws.Range("E21").Value = "  -5.58%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.998"
$ws.Range("E22").Value = "  -0.21%  "
$ws.Range("E23").Value = "  -4.84%  "
$ws.Range("E24").Value = "  -2.62%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.408"
$ws.Range("E25").Value = "  -5.25%  "
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("E27").Value = "  -5.17%  "
$ws.Range("D28").Value = "2.651.74"
$ws.Range("E28").Value = "  -3.14%  "
$ws.Range("D29").Value = "0.0₃0785"
$ws.Range("E29").Value = "  -9.51%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.93"
$ws.Range("E30").Value = "  -6.40%  "
$ws.Range("E31").Value = "  +0.08%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "149.88"
$ws.Range("E32").Value = "  -1.07%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.83"
$ws.Range("E33").Value = "  -5.57%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "18.49"
$ws.Range("E34").Value = "  -5.33%  "
$ws.Range("E35").Value = "  -5.71%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.939"
$ws.Range("E36").Value = "  +5.76%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.91"
$ws.Range("E37").Value = "  -6.77%  "
$ws.Range("E38").Value = "  -7.67%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.03"
$ws.Range("E39").Value = "  -1.17%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.824"
$ws.Range("E40").Value = "  -10.61%  "
$ws.Range("E41").Value = "  -6.93%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "283.24"
$ws.Range("E42").Value = "  -4.01%  "
$ws.Range("E43").Value = "  -7.74%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0995"
$ws.Range("E44").Value = "  -3.27%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.998"
$ws.Range("E45").Value = "  +0.06%  "
$ws.Range("E46").Value = "  -6.25%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0533"
$ws.Range("E47").Value = "  -5.13%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0226"
$ws.Range("E50").Value = "  -5.43%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.54"
$ws.Range("E51").Value = "  -9.82%  "
# Row 48/49 swap: EnergySwap <-> WhiteBITCoin with updated values
$ws.Range("B48").Value = "WhiteBITCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "10.30"
$ws.Range("E48").Value = "  -0.43%  "

$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "18.61"
$ws.Range("E49").Value = "  -5.86%  "
